$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy style from existing header cell (H1) onto the new header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# I / J values per row (row -> [I, J])
$values = @{
    2  = @(1, 5)
    3  = @(1, 5)
    4  = @(1, 6)
    5  = @(1, 6)
    6  = @(1, 7)
    7  = @(1, 7)
    8  = @(1, 7)
    9  = @(1, 6)
    10 = @(1, 7)
    11 = @(1, 6)
    12 = @(1, 7)
    13 = @(1, 6)
    14 = @(1, 7)
    15 = @(1, 6)
    16 = @(1, 6)
    17 = @(1, 7)
    18 = @(1, 7)
    19 = @(1, 6)
    20 = @(1, 5)
    21 = @(1, 9)
    22 = @(1, 4)
    23 = @(1, 6)
    24 = @(1, 6)
    25 = @(1, 8)
    26 = @(1, 7)
    27 = @(1, 6)
    28 = @(1, 5)
    29 = @(1, 5)
    30 = @(1, 5)
    31 = @(5, 8)
    32 = @(4, 7)
    33 = @(5, 8)
    34 = @(4, 6)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
